$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 359.7143
$ws.Range("I80").Value = 99
$ws.Range("K80").Value = 297
$ws.Range("M80").Value = 701
$ws.Range("H83").Value = 359.7143
$ws.Range("I83").Value = 99
$ws.Range("K83").Value = 891
$ws.Range("M83").Value = 4101
$ws.Range("H101").Value = 1291.7693
$ws.Range("I101").Value = 1291.7693
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 3875.3079
$ws.Range("L101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -2253.3079
$ws.Range("H113").Value = 99999.5
$ws.Range("I113").Value = 100000
$ws.Range("K113").Value = 100000
$ws.Range("M113").Value = -96746
$ws.Range("H116").Value = 6917.857
$ws.Range("I116").Value = 6917.5
$ws.Range("J116").Value = 6918
$ws.Range("K116").Value = 6917.5
$ws.Range("L116").Value = 6918
$ws.Range("M116").Value = -3475.5
$ws.Range("N116").Value = -13802
$ws.Range("H137").Value = 2366.1045
$ws.Range("I137").Value = 2359.1428
$ws.Range("J137").Value = 2373.7188
$ws.Range("K137").Value = 7077.428400000001
$ws.Range("L137").Value = 7121.1564
$ws.Range("M137").Value = -4527.428400000001
$ws.Range("N137").Value = -12221.1564
$ws.Range("H138").Value = 3614.1882
$ws.Range("I138").Value = 3589.6875
$ws.Range("K138").Value = 10769.0625
$ws.Range("M138").Value = -5629.0625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 42000
$ws.Range("J7").Value = 42000
$ws.Range("L7").Value = 42000
$ws.Range("N7").Value = -42228
$ws.Range("H32").Value = 15966.652
$ws.Range("I32").Value = 11638.05
$ws.Range("J32").Value = 44824
$ws.Range("K32").Value = 11638.05
$ws.Range("L32").Value = 44824
$ws.Range("M32").Value = -11351.05
$ws.Range("N32").Value = -45398
$ws.Range("H45").Value = 7270.923
$ws.Range("I45").Value = 21319
$ws.Range("J45").Value = 3926.1428
$ws.Range("K45").Value = 21319
$ws.Range("L45").Value = 3926.1428
$ws.Range("M45").Value = -20942
$ws.Range("N45").Value = -4680.1428
$ws.Range("H74").Value = 6614.9414
$ws.Range("I74").Value = 4490.8184
$ws.Range("K74").Value = 4490.8184
$ws.Range("M74").Value = -3616.8184
$ws.Range("H77").Value = 6614.9414
$ws.Range("I77").Value = 4490.8184
$ws.Range("K77").Value = 22454.092
$ws.Range("M77").Value = -18086.092
$ws.Range("H88").Value = 3152.889
$ws.Range("J88").Value = 3111.5715
$ws.Range("L88").Value = 3111.5715
$ws.Range("N88").Value = -3923.5715
$ws.Range("H91").Value = 3152.889
$ws.Range("J91").Value = 3111.5715
$ws.Range("L91").Value = 3111.5715
$ws.Range("N91").Value = -5919.5715
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3365.3333
$ws.Range("I86").Value = 3365.3333
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3365.3333
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -2242.3333
$ws.Range("H89").Value = 3365.3333
$ws.Range("I89").Value = 3365.3333
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 16826.6665
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -11210.6665
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4964.2534
$ws.Range("I31").Value = 4428.2095
$ws.Range("J31").Value = 5787.4644
$ws.Range("K31").Value = 4428.2095
$ws.Range("L31").Value = 5787.4644
$ws.Range("M31").Value = -4133.2095
$ws.Range("N31").Value = -6377.4644
$ws.Range("H34").Value = 4964.2534
$ws.Range("I34").Value = 4428.2095
$ws.Range("J34").Value = 5787.4644
$ws.Range("K34").Value = 4428.2095
$ws.Range("L34").Value = 5787.4644
$ws.Range("M34").Value = -4226.2095
$ws.Range("N34").Value = -6191.4644
$ws.Range("H59").Value = 67499.5
$ws.Range("J59").Value = 70999.39999999999
$ws.Range("L59").Value = 70999.39999999999
$ws.Range("N59").Value = -73289.39999999999
$ws.Range("H60").Value = 29908.6
$ws.Range("J60").Value = 29908.6
$ws.Range("L60").Value = 29908.6
$ws.Range("N60").Value = -30930.6
$ws.Range("H99").Value = 5402.1816
$ws.Range("I99").Value = 7095.4
$ws.Range("J99").Value = 3991.1667
$ws.Range("K99").Value = 7095.4
$ws.Range("L99").Value = 3991.1667
$ws.Range("M99").Value = -5597.4
$ws.Range("N99").Value = -6987.1667
$ws.Range("H122").Value = 1414.6666
$ws.Range("I122").Value = 1397.6
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4192.799999999999
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1742.799999999999
$ws.Range("N122").Value = -9400
$ws.Range("H126").Value = 5402.1816
$ws.Range("I126").Value = 7095.4
$ws.Range("J126").Value = 3991.1667
$ws.Range("K126").Value = 21286.2
$ws.Range("L126").Value = 11973.5001
$ws.Range("M126").Value = -18816.2
$ws.Range("N126").Value = -16913.5001
$ws.Range("H132").Value = 1419.4667
$ws.Range("I132").Value = 1341.0769
$ws.Range("K132").Value = 4023.2307
$ws.Range("M132").Value = -1493.2307
$ws.Range("H134").Value = 2050.05
$ws.Range("I134").Value = 1921.1052
$ws.Range("K134").Value = 5763.3156
$ws.Range("M134").Value = -3228.3156
$ws.Range("H141").Value = 138757.66
$ws.Range("J141").Value = 142792.1
$ws.Range("L141").Value = 142792.1
$ws.Range("N141").Value = -153152.1
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1327.4286
$ws.Range("J5").Value = 1883.6666
$ws.Range("L5").Value = 5650.9998
$ws.Range("N5").Value = -5874.9998
$ws.Range("H46").Value = 28924.514
$ws.Range("I46").Value = 1135.9678
$ws.Range("K46").Value = 3407.9034
$ws.Range("M46").Value = -3316.9034
$ws.Range("H68").Value = 2500.2
$ws.Range("J68").Value = 2500.2
$ws.Range("L68").Value = 7500.599999999999
$ws.Range("N68").Value = -9122.599999999999
$ws.Range("H71").Value = 2500.2
$ws.Range("J71").Value = 2500.2
$ws.Range("L71").Value = 22501.8
$ws.Range("N71").Value = -30613.8
$ws.Range("H107").Value = 1324.1666
$ws.Range("J107").Value = 1361.25
$ws.Range("L107").Value = 4083.75
$ws.Range("N107").Value = -7923.75
$ws.Range("H135").Value = 1327.4286
$ws.Range("J135").Value = 1883.6666
$ws.Range("L135").Value = 16952.9994
$ws.Range("N135").Value = -22022.9994
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 59629.668
$ws.Range("J82").Value = 59949
$ws.Range("L82").Value = 59949
$ws.Range("N82").Value = -60715
$ws.Range("H85").Value = 59629.668
$ws.Range("J85").Value = 59949
$ws.Range("L85").Value = 59949
$ws.Range("N85").Value = -62601
$ws.Range("H117").Value = 12000
$ws.Range("J117").Value = 12000
$ws.Range("L117").Value = 12000
$ws.Range("N117").Value = -18884
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3439.25
$ws.Range("I7").Value = 2751
$ws.Range("K7").Value = 2751
$ws.Range("M7").Value = -2639
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H82").Value = 5696.3335
$ws.Range("I82").Value = 4745.6
$ws.Range("K82").Value = 4745.6
$ws.Range("M82").Value = -4384.6
$ws.Range("H85").Value = 5696.3335
$ws.Range("I85").Value = 4745.6
$ws.Range("K85").Value = 4745.6
$ws.Range("M85").Value = -3497.6
$ws.Range("H122").Value = 3491.8
$ws.Range("I122").Value = 3613.111
$ws.Range("K122").Value = 10839.333
$ws.Range("M122").Value = -8389.332999999999
$ws.Range("H126").Value = 3439.25
$ws.Range("I126").Value = 2751
$ws.Range("K126").Value = 8253
$ws.Range("M126").Value = -5783
$ws.Range("H133").Value = 89996.53
$ws.Range("J133").Value = 89996.53
$ws.Range("L133").Value = 89996.53
$ws.Range("N133").Value = -95056.53
$ws.Range("H136").Value = 3732.3333
$ws.Range("I136").Value = 2459.7334
$ws.Range("K136").Value = 7379.2002
$ws.Range("M136").Value = -4829.2002
$ws.Range("H138").Value = 73140.47
$ws.Range("J138").Value = 81292.53
$ws.Range("L138").Value = 81292.53
$ws.Range("N138").Value = -91572.53
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 69979
$ws.Range("J86").Value = 69979
$ws.Range("L86").Value = 69979
$ws.Range("N86").Value = -72225
$ws.Range("H89").Value = 69979
$ws.Range("J89").Value = 69979
$ws.Range("L89").Value = 349895
$ws.Range("N89").Value = -361127
$ws.Range("H122").Value = 7517.647
$ws.Range("I122").Value = 6838.154
$ws.Range("K122").Value = 20514.462
$ws.Range("M122").Value = -18064.462
$ws.Range("H132").Value = 5076
$ws.Range("I132").Value = 4905.0527
$ws.Range("K132").Value = 14715.1581
$ws.Range("M132").Value = -12185.1581
$ws.Range("H136").Value = 14276.234
$ws.Range("I136").Value = 19984.322
$ws.Range("J136").Value = 3216.8125
$ws.Range("K136").Value = 59952.966
$ws.Range("L136").Value = 9650.4375
$ws.Range("M136").Value = -57402.966
$ws.Range("N136").Value = -14750.4375
